# Generate Report for Handback
# Bump the timestamp strings recorded for the 8cc863a9... handback row:
#   - "Latest HO Xliff Generate Date" (Overview!G2 and de-de!H2 share the
#     same text "2016-08-20 23:08:06" -> "2016-08-20 23:08:56")
#   - zh-cn sheet: Correspond Handoff Datetime (H2) and
#     Correspond Handback DateTime (K2)
#   - de-de sheet: Correspond Handoff Datetime (K2)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date for 8cc863a9-...md
$wsOverview.Range("G2").Value = "2016-08-20 23:08:56"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-20 23:08:51"
$wsZhCn.Range("K2").Value = "2016-08-20 23:09:13"

# de-de: Latest HO Xliff Generate Date (same value as Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-20 23:08:56"

# de-de: Correspond Handoff Datetime
$wsDeDe.Range("K2").Value = "2016-08-20 23:09:20"
